# Portfolio-1 final data update:
# - Regenerate RT (column D) values for all 72 remaining observations
# - Reduce sample size per Block_type x Shift_type cell from n=20 to n=18
#   (drop the last 2 rows of each of the 4 groups -> removes 8 rows total)
# - Renumber PID (column A) 1..18 within each of the 4 groups

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colA = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)
$colB = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2)
$colC = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2)
$colD = @(1124.3317799690401, 1709.4064235687301, 1161.0686566148499, 1265.55928346273, 1371.3963761049199, 1042.64877856463, 1158.8322371244399, 1433.06936614815, 960.36616037058297, 1033.8311251472001, 935.57915299437798, 1210.5998289890799, 1135.0335563932099, 1264.2417308446502, 754.35857150865604, 1120.8502159842999, 1159.7366174062101, 1029.3267165913301, 1665.0309345938902, 2136.7008459000399, 1327.8797183718, 1433.0873272635699, 1573.03722151395, 1345.48615587169, 1236.48009981428, 1612.12934296707, 1245.8357192851902, 1106.1112551853601, 1201.21844325747, 1725.1582336425799, 1372.63625243614, 1530.34167819553, 685.24928887685098, 1333.3773314952898, 1450.24352807265, 1407.22207839672, 1023.94280066857, 1850.7387439409899, 1106.6575050353999, 1120.3937941584099, 1833.87985935918, 1326.0178389372602, 1581.7374897003201, 1504.7991364090501, 1159.73693446109, 1070.7938212614799, 1043.1753931374399, 1063.3954451634299, 1096.77877097294, 1225.0432751395499, 728.25842745163902, 1312.1229895838999, 1055.4949097011399, 1072.96654030129, 1052.0749281753199, 1911.1829429864899, 1293.1088268756901, 1210.82695420966, 1631.0264943521199, 1388.98379007975, 1204.1195660829501, 1554.3478050015201, 1107.4298948481501, 1141.0817476836098, 1194.0285346087301, 1202.65987442761, 1128.48304347558, 1409.6306280656302, 869.45148876735107, 1290.6524151563601, 1476.33555688356, 1214.2383618788299)

# Rows 2..73 hold the 72 surviving observations (18 per Block_type x Shift_type cell)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}

# Drop the now-obsolete trailing rows (old rows 74-81, PID 13-20 of the last group)
$ws.Range("A74:H81").EntireRow.Delete()

# Restore the active-cell selection recorded for this sheet
$ws.Range("E12").Select()

